# "finished parser for local data"
#
# The author finished the local CSV/text-import parser: they took the
# first month's worth of parsed rows from Sheet2 (rows 2-32, i.e. the
# A2:D32 block: Date / Name / Theater(seen-alone "b" columns)) and copied
# them into a brand-new "CSV Export" sheet (A1:D31) placed after Sheet2,
# then left that new sheet active/selected while Sheet2's own selection
# moved to cover that same source block.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(1)          # "Sheet2" - the only sheet so far

# --- add the new "CSV Export" sheet right after Sheet2 ------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$csv = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$csv.Name = "CSV Export"

# --- populate it with the first month of parsed data ---------------------
# Copy formatting first (so column A keeps its date style, B its wrap
# style, C/D their boolean style - even for the one row where Theater /
# Alone are blank), then copy the values/text on top of that.
$src.Range("A2:D32").Copy()
$csv.Range("A1:D31").PasteSpecial(-4122)   # xlPasteFormats

$src.Range("A2:D32").Copy()
$csv.Range("A1:D31").PasteSpecial(-4163)   # xlPasteValues

$excel.CutCopyMode = $false

# --- selections: Sheet2 now shows the block that was copied out ---------
[void]$src.Activate()
$src.Range("A2:D32").Select()

# --- the new sheet ends up the active / front-most tab -------------------
[void]$csv.Activate()
$csv.Range("A1:D31").Select()
